# Weekly data entry: append the row for 2025-04-24 (serial 45771) on each
# of the four tracking sheets, then leave the selection/active sheet the
# way the author left it (ends on "Saldo", cell B19 selected).

$wb = $excel.ActiveWorkbook

# --- Produzione ---------------------------------------------------------
$ws = $wb.Worksheets.Item("Produzione")
$ws.Activate() | Out-Null
$ws.Range("A18").Value = 45771
$ws.Range("B18").NumberFormat = $ws.Range("B17").NumberFormat
$ws.Range("B18").Value = 26458.32
$ws.Range("A19").Select() | Out-Null

# --- Entrate -------------------------------------------------------------
$ws = $wb.Worksheets.Item("Entrate")
$ws.Activate() | Out-Null
$ws.Range("A18").Value = 45771
$ws.Range("B18").NumberFormat = $ws.Range("B17").NumberFormat
$ws.Range("B18").Value = 17728.89
$ws.Range("B19").Select() | Out-Null

# --- Uscite ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("Uscite")
$ws.Activate() | Out-Null
$ws.Range("A18").Value = 45771
$ws.Range("B18").NumberFormat = $ws.Range("B17").NumberFormat
$ws.Range("B18").Value = 3343.54
$ws.Range("C21").Select() | Out-Null

# --- Saldo ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Saldo")
$ws.Activate() | Out-Null
$ws.Range("A18").Value = 45771
$ws.Range("B18").NumberFormat = $ws.Range("B17").NumberFormat
$ws.Range("B18").Value = 56863.88
$ws.Range("B19").Select() | Out-Null
